$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header-like numeric values
$ws.Range("B2").Value = 1
$ws.Range("E2").Value = 2
$ws.Range("H2").Value = 0

# Row 4 data values
$ws.Range("B4").Value = 0.5607188170993559
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.1214376341987118
$ws.Range("E4").Value = 0.7831777282599498
$ws.Range("G4").Value = 0.5663554565198996
$ws.Range("H4").Value = 0.7087799791449426
$ws.Range("I4").Value = -1
$ws.Range("J4").Value = 0.4175599582898852
